$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = "2018-12-31 00:00:00"

$ws.Range("O2").Value = 228378431.75
$ws.Range("P2").Value = 733579538.71
$ws.Range("Q2").Value = 443958944.34
$ws.Range("R2").Value = 20.6133084115
$ws.Range("S2").Value = 296727498.68
$ws.Range("T2").Value = 296727498.68
$ws.Range("U2").Value = 23.915278286
$ws.Range("V2").Value = 95178485.92
$ws.Range("W2").Value = 38390689.43
$ws.Range("X2").Value = 4618982.87
$ws.Range("Y2").Value = 292731929.84
$ws.Range("Z2").Value = 284793463.36
$ws.Range("AA2").Value = 56415031.61
$ws.Range("AG2").Value = 9043287.439999999
$ws.Range("AP2").Value = 38.8187918623
$ws.Range("AQ2").Value = 88.619165177805
$ws.Range("AR2").Value = 93.02671005351399
$ws.Range("AS2").Value = 230162786.96
$ws.Range("AT2").Value = 97.34895795466799
